# Update forecast summary: fill in the previously-empty "Amazon Mean Forecast"
# column (D) and refresh the P70/P80/P90 forecast columns (E/F/G) now that the
# Auto Arima model has been removed from the pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$data = @{
    2  = @(75, 88, 101, 121)
    3  = @(65, 78, 91, 111)
    4  = @(66, 79, 92, 113)
    5  = @(65, 78, 91, 112)
    6  = @(67, 81, 96, 121)
    7  = @(67, 81, 95, 118)
    8  = @(69, 83, 99, 124)
    9  = @(71, 86, 103, 130)
    10 = @(69, 83, 99, 124)
    11 = @(70, 85, 101, 128)
    12 = @(71, 87, 104, 133)
    13 = @(75, 91, 111, 142)
    14 = @(72, 88, 106, 135)
    15 = @(72, 88, 108, 140)
    16 = @(70, 85, 104, 135)
    17 = @(68, 83, 101, 131)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 4).Value = $values[0]  # D: Amazon Mean Forecast
    $ws.Cells.Item($row, 5).Value = $values[1]  # E: Amazon P70 Forecast
    $ws.Cells.Item($row, 6).Value = $values[2]  # F: Amazon P80 Forecast
    $ws.Cells.Item($row, 7).Value = $values[3]  # G: Amazon P90 Forecast
}
